$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-03-05 Wednesday"; new="2025-03-06 Thursday"},
    @{old="390÷2="; new="559÷7="},
    @{old="505÷8="; new="401÷9="},
    @{old="675÷7="; new="457÷9="},
    @{old="835÷9="; new="377÷5="},
    @{old="845÷7="; new="141÷5="},
    @{old="835÷5="; new="196÷9="},
    @{old="670÷6="; new="821÷8="},
    @{old="678÷5="; new="889÷9="},
    @{old="727÷8="; new="983÷3="},
    @{old="531÷9="; new="725÷5="},
    @{old="789÷7="; new="141÷5="},
    @{old="267÷4="; new="259÷4="},
    @{old="786÷8="; new="464÷9="},
    @{old="177÷9="; new="964÷2="},
    @{old="840÷2="; new="532÷6="},
    @{old="180÷7="; new="514÷7="},
    @{old="135÷8="; new="854÷2="},
    @{old="945÷7="; new="669÷8="},
    @{old="434÷2="; new="130÷5="},
    @{old="941÷7="; new="501÷9="},
    @{old="414÷6="; new="766÷2="},
    @{old="795÷5="; new="119÷9="},
    @{old="349÷5="; new="893÷7="},
    @{old="867÷4="; new="809÷6="},
    @{old="163÷9="; new="639÷8="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
